$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 517
$ws.Range("I2").Value = 410
$ws.Range("J2").Value = 945
$ws.Range("K2").Value = 410
$ws.Range("L2").Value = 945
$ws.Range("M2").Value = -297
$ws.Range("N2").Value = -1171

$ws.Range("H6").Value = 532.9167
$ws.Range("I6").Value = 4.3333335
$ws.Range("J6").Value = 1061.5
$ws.Range("K6").Value = 13.0000005
$ws.Range("L6").Value = 3184.5
$ws.Range("M6").Value = 98.9999995
$ws.Range("N6").Value = -3408.5

$ws.Range("H9").Value = 195.33333
$ws.Range("I9").Value = 193
$ws.Range("J9").Value = 200
$ws.Range("K9").Value = 193
$ws.Range("L9").Value = 200
$ws.Range("M9").Value = -24
$ws.Range("N9").Value = -538

$ws.Range("H21").Value = 9999.5
$ws.Range("I21").Value = 9999.5
$ws.Range("K21").Value = 9999.5
$ws.Range("M21").Value = -9531.5

$ws.Range("H23").Value = 9999.5
$ws.Range("I23").Value = 9999.5
$ws.Range("K23").Value = 9999.5
$ws.Range("M23").Value = -9765.5

$ws.Range("H31").Value = 11.333333
$ws.Range("I31").Value = 11.333333
$ws.Range("K31").Value = 33.999999
$ws.Range("M31").Value = 196.000001

$ws.Range("H38").Value = 1668
$ws.Range("I38").Value = 311.42856
$ws.Range("K38").Value = 934.28568
$ws.Range("M38").Value = -562.28568

$ws.Range("H39").Value = 405.18182
$ws.Range("I39").Value = 222.71428
$ws.Range("K39").Value = 668.14284
$ws.Range("M39").Value = -372.14284

$ws.Range("H58").Value = 1448.5555
$ws.Range("J58").Value = 3187.5
$ws.Range("L58").Value = 9562.5
$ws.Range("N58").Value = -9862.5

$ws.Range("H80").Value = 568.2222
$ws.Range("J80").Value = 567.9
$ws.Range("L80").Value = 1703.7
$ws.Range("N80").Value = -3699.7

$ws.Range("H83").Value = 568.2222
$ws.Range("J83").Value = 567.9
$ws.Range("L83").Value = 5111.099999999999
$ws.Range("N83").Value = -15095.1

$ws.Range("H87").Value = 91428.60000000001
$ws.Range("J87").Value = 91428.60000000001
$ws.Range("L87").Value = 91428.60000000001
$ws.Range("N87").Value = -93924.60000000001

$ws.Range("H90").Value = 91428.60000000001
$ws.Range("J90").Value = 91428.60000000001
$ws.Range("L90").Value = 274285.8
$ws.Range("N90").Value = -286765.8

$ws.Range("H112").Value = 3459.3333
$ws.Range("J112").Value = 3459.3333
$ws.Range("L112").Value = 10377.9999
$ws.Range("N112").Value = -12593.9999

$ws.Range("H130").Value = 110000
$ws.Range("J130").Value = 110000
$ws.Range("L130").Value = 110000
$ws.Range("N130").Value = -120040

$ws.Range("H132").Value = 3969.88
$ws.Range("I132").Value = 1147.591
$ws.Range("K132").Value = 3442.773
$ws.Range("M132").Value = -912.7729999999997

$ws.Range("H138").Value = 2802.95
$ws.Range("I138").Value = 1733.9
$ws.Range("J138").Value = 3872
$ws.Range("K138").Value = 5201.700000000001
$ws.Range("L138").Value = 11616
$ws.Range("M138").Value = -61.70000000000073
$ws.Range("N138").Value = -21896

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 4500
$ws.Range("I61").Value = 4500
$ws.Range("K61").Value = 4500
$ws.Range("M61").Value = -4288

$ws.Range("H122").Value = 1585.2354
$ws.Range("I122").Value = 1449.9375
$ws.Range("K122").Value = 4349.8125
$ws.Range("M122").Value = -1899.8125

$ws.Range("H136").Value = 4500
$ws.Range("I136").Value = 4500
$ws.Range("K136").Value = 13500
$ws.Range("M136").Value = -10950

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("I20").Value = 10500
$ws.Range("J20").Value = 1000
$ws.Range("K20").Value = 10500
$ws.Range("L20").Value = 1000
$ws.Range("M20").Value = -10253
$ws.Range("N20").Value = -1494

$ws.Range("H86").Value = 3422.125
$ws.Range("I86").Value = 1350.1818
$ws.Range("K86").Value = 1350.1818
$ws.Range("M86").Value = -227.1818000000001

$ws.Range("H89").Value = 3422.125
$ws.Range("I89").Value = 1350.1818
$ws.Range("K89").Value = 6750.909000000001
$ws.Range("M89").Value = -1134.909000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H5").Value = 893.4286
$ws.Range("I5").Value = 217.66667
$ws.Range("J5").Value = 2109.8
$ws.Range("K5").Value = 217.66667
$ws.Range("L5").Value = 2109.8
$ws.Range("M5").Value = -105.66667
$ws.Range("N5").Value = -2333.8

$ws.Range("H9").Value = 175000
$ws.Range("J9").Value = 175000
$ws.Range("L9").Value = 175000
$ws.Range("N9").Value = -175336

$ws.Range("H23").Value = 9
$ws.Range("I23").Value = 9
$ws.Range("J23").Value = 0
$ws.Range("K23").Value = 9
$ws.Range("L23").Value = 0
$ws.Range("M23").ClearContents()
$ws.Range("N23").Value = 231

$ws.Range("H25").Value = 3381.8572
$ws.Range("I25").Value = 3149.6667
$ws.Range("J25").Value = 4775
$ws.Range("K25").Value = 3149.6667
$ws.Range("L25").Value = 4775
$ws.Range("M25").Value = -2975.6667
$ws.Range("N25").Value = -5123

$ws.Range("H27").Value = 9
$ws.Range("I27").Value = 9
$ws.Range("J27").Value = 0
$ws.Range("K27").Value = 9
$ws.Range("L27").Value = 0
$ws.Range("M27").ClearContents()
$ws.Range("N27").Value = 183

$ws.Range("H62").Value = 2277.75
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()

$ws.Range("H65").Value = 2277.75
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()

$ws.Range("H86").Value = 3999.5
$ws.Range("I86").Value = 3999
$ws.Range("K86").Value = 3999
$ws.Range("M86").Value = -2876

$ws.Range("H89").Value = 3999.5
$ws.Range("I89").Value = 3999
$ws.Range("K89").Value = 19995
$ws.Range("M89").Value = -14379

$ws.Range("H107").Value = 843.1539
$ws.Range("I107").Value = 752.3333
$ws.Range("J107").Value = 1047.5
$ws.Range("K107").Value = 752.3333
$ws.Range("L107").Value = 1047.5
$ws.Range("M107").Value = 1167.6667
$ws.Range("N107").Value = -4887.5

$ws.Range("H132").Value = 3859.7058
$ws.Range("I132").Value = 3593.9285
$ws.Range("K132").Value = 10781.7855
$ws.Range("M132").Value = -8251.7855

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 19.571428
$ws.Range("I2").Value = 22.714285
$ws.Range("J2").Value = 16.428572
$ws.Range("K2").Value = 136.28571
$ws.Range("L2").Value = 98.57143199999999
$ws.Range("M2").Value = -23.28570999999999
$ws.Range("N2").Value = -324.571432

$ws.Range("H4").Value = 284127
$ws.Range("I4").Value = 284127
$ws.Range("K4").Value = 852381
$ws.Range("M4").Value = -852269

$ws.Range("H131").Value = 2324.0625
$ws.Range("I131").Value = 2096.5
$ws.Range("J131").Value = 2460.6
$ws.Range("K131").Value = 6289.5
$ws.Range("L131").Value = 7381.799999999999
$ws.Range("M131").Value = -1249.5
$ws.Range("N131").Value = -17461.8

$ws.Range("H140").Value = 2003.75
$ws.Range("I140").Value = 1718.5714
$ws.Range("J140").Value = 4000
$ws.Range("K140").Value = 5155.7142
$ws.Range("L140").Value = 12000
$ws.Range("M140").Value = 24.28579999999965
$ws.Range("N140").Value = -22360

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 0
$ws.Range("I70").Value = 0
$ws.Range("K70").Value = 0
$ws.Range("M70").ClearContents()

$ws.Range("H73").Value = 0
$ws.Range("I73").Value = 0
$ws.Range("K73").Value = 0
$ws.Range("M73").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 6400
$ws.Range("I16").Value = 6500
$ws.Range("J16").Value = 6000
$ws.Range("K16").Value = 6500
$ws.Range("L16").Value = 6000
$ws.Range("M16").Value = -6330
$ws.Range("N16").Value = -6340

$ws.Range("H46").Value = 2710.818
$ws.Range("J46").Value = 3217.1333
$ws.Range("L46").Value = 3217.1333
$ws.Range("N46").Value = -3593.1333

$ws.Range("H136").Value = 3100
$ws.Range("I136").Value = 2400
$ws.Range("K136").Value = 7200
$ws.Range("M136").Value = -4650

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 4413.5
$ws.Range("I126").Value = 2264.5
$ws.Range("J126").Value = 6562.5
$ws.Range("K126").Value = 6793.5
$ws.Range("L126").Value = 19687.5
$ws.Range("M126").Value = -4323.5
$ws.Range("N126").Value = -24627.5
